$d = $word.ActiveDocument
$d.Content.Find.Execute("000", $true, $false, $false, $false, $false, $true, 1, $false, "{{anuncio}}", 2)
